$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-22 01:49:36"

# Drop every existing hyperlink up front. (In this engine, deleting via a
# range's Hyperlinks collection clears the whole sheet's collection rather
# than just that cell's link, and row-insert doesn't relocate hyperlink
# anchors -- so the only clean way to end up with correctly targeted links
# after the insert is to rebuild the full set once the new row layout and
# URL text are in place.)
$ws.Range("A1").Hyperlinks.Delete()

# Insert a new row at position 10, shifting existing rows 10-12 down to 11-13
$ws.Rows.Item(10).Insert()

# Update timestamps for all data rows (2-13) to the new timestamp
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# Fill in the newly inserted row 10 with its data
$ws.Cells.Item(10, 2).Value = "Excelやスプレッドシートでのデータシュミレーション クエリ(query)や関数利用"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5397980"
$ws.Cells.Item(10, 7).Value = 10

# Re-create the hyperlinks for F2:F13 in order so relationship ids come out
# sequential (rId1..rId12) and every link's target matches the URL text now
# shown in its row.
$urls = @(
    "https://www.lancers.jp/work/detail/5397594",
    "https://www.lancers.jp/work/detail/5397680",
    "https://www.lancers.jp/work/detail/5397543",
    "https://www.lancers.jp/work/detail/5397930",
    "https://www.lancers.jp/work/detail/5397812",
    "https://www.lancers.jp/work/detail/5397452",
    "https://www.lancers.jp/work/detail/5397615",
    "https://www.lancers.jp/work/detail/5397887",
    "https://www.lancers.jp/work/detail/5397980",
    "https://www.lancers.jp/work/detail/5397962",
    "https://www.lancers.jp/work/detail/5397817",
    "https://www.lancers.jp/work/detail/5395809"
)
for ($i = 0; $i -lt $urls.Length; $i++) {
    $row = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $urls[$i])
}

# Widen column B from 40 to 45 (character units). The engine quantizes
# ColumnWidth to whole pixels on save, so 45.0 round-trips to 45.8333; a
# value a bit under 45 lands on the same pixel bucket as a clean 45.
$ws.Columns.Item(2).ColumnWidth = 44.15
